# Daily attendance processing - 2025-10-10 11:18:31
# Swap the order of the comma-separated "Recorded By" values in column G
# (the last-listed recorder moves to the front of the list) for the
# affected rows on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = 3,6,10,11,12,13,14,15,30,33,37,38,39,40,41,42,57,60,64,65,66,67,68,69,86,87,88,89,90,93,95,112,113,114,115,116,119,121,138,139,140,141,142,145,147

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $current = $cell.Value()
    $parts = $current.Split(",")
    $first = $parts[0].Trim()
    $second = $parts[1].Trim()
    $cell.Value = "$second, $first"
}
